$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$values = @(
    "94-25=", "53-9=", "39+54=", "74+17=", "84-68=", "54-39=", "91-12=", "80-31=", "91-54=", "80-43=",
    "46-37=", "82-69=", "44-8=", "4+19=", "29+42=", "90-13=", "92-45=", "9+78=", "9+24=", "36+27=",
    "86-7=", "63-57=", "37+39=", "50-35=", "60-55=", "84-28=", "61-6=", "80-2=", "75-36=", "65-37=",
    "31-6=", "63-6=", "62-13=", "16+56=", "49+32=", "70-48=", "16+25=", "25-17=", "30-9=", "17+6=",
    "28+26=", "43-6=", "55+38=", "70-58=", "34+29=", "95-67=", "61-44=", "70-43=", "89+4=", "60-31=",
    "96-18=", "87+9=", "31-13=", "79+15=", "9+55=", "90-7=", "78+8=", "63-19=", "82-79=", "37+25=",
    "58+36=", "92-15=", "43+19=", "37+45=", "70-62=", "42-4=", "46-17=", "3+59=", "12+19=", "29+65=",
    "7+44=", "94-6=", "74-49=", "59+38=", "48+19=", "46-8=", "80-7=", "29+54=", "11-2=", "43-7=",
    "88+4=", "55+6=", "66-19=", "41-22=", "48+24=", "61-5=", "27+26=", "22-18=", "38+39=", "9+52=",
    "73-6=", "39+47=", "59+27=", "13-5=", "59+19=", "9+44=", "58+29=", "19+13=", "8+14=", "34+28="
)

$idx = 0
for ($r = 1; $r -le 20; $r++) {
    for ($c = 1; $c -le 5; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $values[$idx]
        $idx++
    }
}

Write-Host "Done. idx=" $idx
